$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM")
$ws.Range("D17").Value = "Line1`nLine2`nLine3"
$v = $ws.Range("D17").Value2
Write-Host "D17:" $v
